$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd14"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 54.82987066666667
$ws.Range("H2").Value = 164.489612
$ws.Range("I2").Value = 0.9762630652055824
$ws.Range("J2").Value = 0.9762630652055824
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.906497
$ws.Range("N2").Value = 68.719491
$ws.Range("O2").Value = 0.9446038650914245
$ws.Range("P2").Value = 0.9446038650914245
$ws.Range("Q2").Value = 1255.960267936388
$ws.Range("R2").Value = 11303.64241142749
$ws.Range("S2").Value = 0.9221818647391945
$ws.Range("T2").Value = 0.9221818647391945

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd14"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 54.82987066666667
$ws.Range("H3").Value = 164.489612
$ws.Range("I3").Value = 0.9762630652055824
$ws.Range("J3").Value = 0.9762630652055824
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1329193333333333
$ws.Range("N3").Value = 0.3987579999999999
$ws.Range("O3").Value = 0.005481244732096839
$ws.Range("P3").Value = 0.005481244732096839
$ws.Range("Q3").Value = 7.287949855766222
$ws.Range("R3").Value = 65.591548701896
$ws.Range("S3").Value = 0.005351136783298811
$ws.Range("T3").Value = 0.005351136783298811

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cd14"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 54.82987066666667
$ws.Range("H4").Value = 164.489612
$ws.Range("I4").Value = 0.9762630652055824
$ws.Range("J4").Value = 0.9762630652055824
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.210428333333333
$ws.Range("N4").Value = 3.631285
$ws.Range("O4").Value = 0.04991489017647865
$ws.Range("P4").Value = 0.04991489017647865
$ws.Range("Q4").Value = 66.36762896793557
$ws.Range("R4").Value = 597.3086607114201
$ws.Range("S4").Value = 0.04873006368308906
$ws.Range("T4").Value = 0.04873006368308905

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd14"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9891043333333333
$ws.Range("H5").Value = 2.967313
$ws.Range("I5").Value = 0.01761131325912771
$ws.Range("J5").Value = 0.01761131325912771
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.906497
$ws.Range("N5").Value = 68.719491
$ws.Range("O5").Value = 0.9446038650914245
$ws.Range("P5").Value = 0.9446038650914245
$ws.Range("Q5").Value = 22.656915444187
$ws.Range("R5").Value = 203.912238997683
$ws.Range("S5").Value = 0.01663571457390788
$ws.Range("T5").Value = 0.01663571457390788

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cd14"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9891043333333333
$ws.Range("H6").Value = 2.967313
$ws.Range("I6").Value = 0.01761131325912771
$ws.Range("J6").Value = 0.01761131325912771
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1329193333333333
$ws.Range("N6").Value = 0.3987579999999999
$ws.Range("O6").Value = 0.005481244732096839
$ws.Range("P6").Value = 0.005481244732096839
$ws.Range("Q6").Value = 0.1314710885837777
$ws.Range("R6").Value = 1.183239797254
$ws.Range("S6").Value = 0.00009653191802690096
$ws.Range("T6").Value = 0.00009653191802690096

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cd14"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9891043333333333
$ws.Range("H7").Value = 2.967313
$ws.Range("I7").Value = 0.01761131325912771
$ws.Range("J7").Value = 0.01761131325912771
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.210428333333333
$ws.Range("N7").Value = 3.631285
$ws.Range("O7").Value = 0.04991489017647865
$ws.Range("P7").Value = 0.04991489017647865
$ws.Range("Q7").Value = 1.197239909689444
$ws.Range("R7").Value = 10.775159187205
$ws.Range("S7").Value = 0.0008790667671929218
$ws.Range("T7").Value = 0.0008790667671929217

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd14"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3440333333333334
$ws.Range("H8").Value = 1.0321
$ws.Range("I8").Value = 0.00612562153528991
$ws.Range("J8").Value = 0.006125621535289909
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.906497
$ws.Range("N8").Value = 68.719491
$ws.Range("O8").Value = 0.9446038650914245
$ws.Range("P8").Value = 0.9446038650914245
$ws.Range("Q8").Value = 7.880598517900001
$ws.Range("R8").Value = 70.92538666110001
$ws.Range("S8").Value = 0.005786285778322115
$ws.Range("T8").Value = 0.005786285778322114

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd14"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3440333333333334
$ws.Range("H9").Value = 1.0321
$ws.Range("I9").Value = 0.00612562153528991
$ws.Range("J9").Value = 0.006125621535289909
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1329193333333333
$ws.Range("N9").Value = 0.3987579999999999
$ws.Range("O9").Value = 0.005481244732096839
$ws.Range("P9").Value = 0.005481244732096839
$ws.Range("Q9").Value = 0.0457286813111111
$ws.Range("R9").Value = 0.4115581317999999
$ws.Range("S9").Value = 0.00003357603077112677
$ws.Range("T9").Value = 0.00003357603077112677

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cd14"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3440333333333334
$ws.Range("H10").Value = 1.0321
$ws.Range("I10").Value = 0.00612562153528991
$ws.Range("J10").Value = 0.006125621535289909
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.210428333333333
$ws.Range("N10").Value = 3.631285
$ws.Range("O10").Value = 0.04991489017647865
$ws.Range("P10").Value = 0.04991489017647865
$ws.Range("Q10").Value = 0.4164276942777779
$ws.Range("R10").Value = 3.7478492485
$ws.Range("S10").Value = 0.0003057597261966684
$ws.Range("T10").Value = 0.0003057597261966683

